$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.733.09'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.86%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.847.43'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.30%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.009'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '336.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.53%  '

$ws.Range("E6").Value = '  +0.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4658'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.84%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3857'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.72'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.76%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07897'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.06%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9669'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.98%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.25'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.10%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.846.49'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.884'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.31%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.139'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.15%  '

$ws.Range("E16").Value = '  +0.11%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.02'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.12%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06623'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.59%  '

$ws.Range("E19").Value = '  -0.46%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.009'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.25%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '27.727.23'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.84%  '

$ws.Range("E23").Value = '  -0.81%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.81'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.56%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.295'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.60%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.070.50'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.32%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '158.86'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.00%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.45'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.17%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.060'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.16%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.388'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.25%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '118.59'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.17%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09440'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.51%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9402'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.603'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.24%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.250'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.86%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.325'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.35%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06024'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.24%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02202'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.96%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.224'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.45%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.008'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.152'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.33%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5800'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.61%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1843'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.82%  '

$ws.Range("E44").Value = '  -2.81%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.283'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.60%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5443'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.39%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.91'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.67%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.937'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.64%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06829'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.07%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '110.74'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.70%  '

$ws.Range("E51").Value = '  -32.12%  '
